$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("testSheet")
$ws2 = $wb.Worksheets.Item("paymentRequest")

# Remove 3 now-unused blank rows between the "delete" sample header (row 16)
# and the "delete" statement block (old row 21), collapsing the gap from
# rows 17-20 down to just row 17. Everything below shifts up by 3 rows.
$ws1.Rows("18:20").Delete()

# The old "sequence-restart" block is now the "alter-sequence" block - rename
# the label (a brand-new shared string; the original "sequence-restart" text
# is still used by the paymentRequest sheet).
$ws1.Range("A23").Value = "alter-sequence"

# Bold the three sample "insert" data rows.
$ws1.Range("A3:K3").Font.Bold = $true
$ws1.Range("A8:K8").Font.Bold = $true
$ws1.Range("A13:K13").Font.Bold = $true

# Bold the entity-name and id cells of the two statement blocks.
$ws1.Range("A20").Font.Bold = $true
$ws1.Range("A25").Font.Bold = $true

# Zoom both sheets in a bit.
$ws2.Activate()
$excel.ActiveWindow.Zoom = 130

$ws1.Activate()
$excel.ActiveWindow.Zoom = 130
$ws1.Range("C26").Select() | Out-Null
